$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) ${Company} paragraph: font size 40 half-pts (20pt) -> 36 half-pts (18pt)
#    (covers both the run and the paragraph mark / pPr rPr)
# ------------------------------------------------------------------
$rngCompany = $d.Content
$rngCompany.Find.ClearFormatting()
if ($rngCompany.Find.Execute("`${Company}")) {
    $paraCompany = $rngCompany.Paragraphs.Item(1).Range
    $paraCompany.Font.Size = 18
}

# ------------------------------------------------------------------
# 2) ${PCode} - ${Project} runs: font size 36 half-pts (18pt) -> 32 half-pts (16pt)
#    Also wrap the two runs in a "_GoBack" bookmark.
# ------------------------------------------------------------------
$rngPCode = $d.Content
$rngPCode.Find.ClearFormatting()
[void]$rngPCode.Find.Execute("`${PCode} - ")
$startMark = $rngPCode.Start

$rngProject = $d.Content
$rngProject.Find.ClearFormatting()
[void]$rngProject.Find.Execute("`${Project}")
$endMark = $rngProject.End

$rngPCodeProject = $d.Range($startMark, $endMark)
$rngPCodeProject.Font.Size = 16

# ------------------------------------------------------------------
# 3) ${Invoice} paragraph: font size 40 half-pts (20pt) -> 36 half-pts (18pt)
# ------------------------------------------------------------------
$rngInvoice = $d.Content
$rngInvoice.Find.ClearFormatting()
if ($rngInvoice.Find.Execute("`${Invoice}")) {
    $paraInvoice = $rngInvoice.Paragraphs.Item(1).Range
    $paraInvoice.Font.Size = 18
}

# ------------------------------------------------------------------
# 4) Remove the old "_GoBack" bookmark (was sitting just before the
#    "${S}     s/d     ${E}" run) and re-add it around the
#    "${PCode} - ${Project}" range.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
[void]$d.Bookmarks.Add("_GoBack", $rngPCodeProject)
